$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add formulas to columns G, H, I for rows 2-4
$ws.Range("G2").Formula = "=ROW()"
$ws.Range("H2").Formula = "=SUM(B2,F2)"
$ws.Range("I2").Formula = "=SUM(B`$4,F`$4)"

$ws.Range("G3").Formula = "=ROW()"
$ws.Range("H3").Formula = "=SUM(B3,F3)"
$ws.Range("I3").Formula = "=SUM(B`$4,F`$4)"

$ws.Range("G4").Formula = "=ROW()"
$ws.Range("H4").Formula = "=SUM(B4,F4)"
$ws.Range("I4").Formula = "=SUM(B`$4,F`$4)"

$wb.Save()
